$d = $word.ActiveDocument

# 1) "modyfikować " and "Grupy kursów " were two separate runs with identical
#    bold formatting; collapse them into a single run with the combined text
#    by doing a same-formatting find & replace over the exact span.
$d.Content.Find.Execute("modyfikować Grupy kursów ", $true, $false, $false, $false, $false, $true, 1, $false, "modyfikować Grupy kursów ", 2) | Out-Null

# 2) The built-in "Normal" style should no longer auto-hyphenate -- turn off
#    paragraph hyphenation, which serializes as <w:suppressAutoHyphens/>.
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Hyphenation = $false
